$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"). Copy the existing header
# style from H1 (bold / bordered / centered) so the new headers match
# the look of the rest of the header row, then overwrite the text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# New data columns I ("I0") and J ("IF") for rows 2-65.
$iVals = @(5,8,8,8,7,8,7,7,7,7,6,7,6,8,8,7,7,6,7,7,8,7,6,6,8,10,7,7,7,8,8,7,7,11,7,5,6,6,7,6,7,7,7,7,6,8,7,7,9,7,7,7,6,6,8,9,6,7,6,6,4,7,7,3)
$jVals = @(5,8,8,8,7,8,7,7,7,7,6,7,6,8,8,7,7,6,7,7,8,7,7,6,8,10,7,7,7,8,8,7,7,11,7,5,6,6,7,7,7,7,7,7,6,8,7,7,9,7,7,7,6,6,8,9,6,7,6,6,4,7,7,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
